$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: force a cell to keep an exact literal text value
# (Excel auto-converts numeric-looking strings to numbers via .Value,
#  which rounds/reformats them - e.g. "9.000" -> 9, "1.001" -> 1.0009999999999999).
# Setting NumberFormat to Text ("@") first forces literal storage, then we
# restore the "Normal" style so no stray style index is left on the cell.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.437.03"
Set-TextValue $ws.Range("E2") "  +0.71%  "

Set-TextValue $ws.Range("D3") "1.874.21"
Set-TextValue $ws.Range("E3") "  +0.82%  "

Set-TextValue $ws.Range("E4") "  -0.02%  "

Set-TextValue $ws.Range("D5") "0.7179"
Set-TextValue $ws.Range("E5") "  +0.86%  "

Set-TextValue $ws.Range("D6") "239.39"
Set-TextValue $ws.Range("E6") "  +0.59%  "

Set-TextValue $ws.Range("D7") "1.001"
Set-TextValue $ws.Range("E7") "  +0.02%  "

Set-TextValue $ws.Range("D8") "0.07821"
Set-TextValue $ws.Range("E8") "  -3.55%  "

Set-TextValue $ws.Range("D9") "0.3075"
Set-TextValue $ws.Range("E9") "  +1.14%  "

Set-TextValue $ws.Range("D10") "25.37"
Set-TextValue $ws.Range("E10") "  +9.39%  "

Set-TextValue $ws.Range("E11") "  +0.60%  "

Set-TextValue $ws.Range("D12") "1.874.79"
Set-TextValue $ws.Range("E12") "  +1.03%  "

$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D13") "0.7237"
Set-TextValue $ws.Range("E13") "  +2.30%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D14") "5.243"
Set-TextValue $ws.Range("E14") "  +1.37%  "

Set-TextValue $ws.Range("D15") "90.72"
Set-TextValue $ws.Range("E15") "  +1.30%  "

Set-TextValue $ws.Range("D16") "29.500.83"
Set-TextValue $ws.Range("E16") "  +0.89%  "

Set-TextValue $ws.Range("D17") "5.857"
Set-TextValue $ws.Range("E17") "  +1.11%  "

Set-TextValue $ws.Range("D18") "0.000007866"
Set-TextValue $ws.Range("E18") "  -0.31%  "

Set-TextValue $ws.Range("D19") "241.97"
Set-TextValue $ws.Range("E19") "  +1.98%  "

Set-TextValue $ws.Range("E20") "  -0.50%  "

Set-TextValue $ws.Range("D21") "2.132.15"
Set-TextValue $ws.Range("E21") "  +0.74%  "

Set-TextValue $ws.Range("D22") "0.9998"
Set-TextValue $ws.Range("E22") "  -0.10%  "

Set-TextValue $ws.Range("E23") "  -0.08%  "

Set-TextValue $ws.Range("D24") "7.758"
Set-TextValue $ws.Range("E24") "  +4.48%  "

Set-TextValue $ws.Range("D25") "0.1559"
Set-TextValue $ws.Range("E25") "  +6.65%  "

Set-TextValue $ws.Range("D26") "163.16"
Set-TextValue $ws.Range("E26") "  +0.39%  "

Set-TextValue $ws.Range("D27") "9.000"
Set-TextValue $ws.Range("E27") "  +0.36%  "

Set-TextValue $ws.Range("D28") "18.33"
Set-TextValue $ws.Range("E28") "  +1.27%  "

Set-TextValue $ws.Range("D29") "1.935"
Set-TextValue $ws.Range("E29") "  -1.11%  "

Set-TextValue $ws.Range("D30") "1.358"
Set-TextValue $ws.Range("E30") "  -4.97%  "

Set-TextValue $ws.Range("E31") "  +0.03%  "

Set-TextValue $ws.Range("D32") "4.333"
Set-TextValue $ws.Range("E32") "  -1.59%  "

Set-TextValue $ws.Range("E33") "  +1.83%  "

Set-TextValue $ws.Range("D34") "0.05257"
Set-TextValue $ws.Range("E34") "  +0.67%  "

Set-TextValue $ws.Range("D35") "1.200"
Set-TextValue $ws.Range("E35") "  +2.73%  "

Set-TextValue $ws.Range("D36") "0.7184"
Set-TextValue $ws.Range("E36") "  +1.42%  "

Set-TextValue $ws.Range("D37") "1.005"
Set-TextValue $ws.Range("E37") "  +0.56%  "

Set-TextValue $ws.Range("D38") "2.675"
Set-TextValue $ws.Range("E38") "  +0.08%  "

Set-TextValue $ws.Range("D39") "0.01867"
Set-TextValue $ws.Range("E39") "  +0.40%  "

Set-TextValue $ws.Range("E40") "  -0.29%  "

Set-TextValue $ws.Range("D41") "1.182.49"
Set-TextValue $ws.Range("E41") "  +3.57%  "

Set-TextValue $ws.Range("E42") "  -1.53%  "

Set-TextValue $ws.Range("D43") "72.40"
Set-TextValue $ws.Range("E43") "  +3.21%  "

Set-TextValue $ws.Range("D44") "6.015"
Set-TextValue $ws.Range("E44") "  +2.47%  "

Set-TextValue $ws.Range("D45") "0.4315"
Set-TextValue $ws.Range("E45") "  +0.71%  "

Set-TextValue $ws.Range("D46") "1.001"
Set-TextValue $ws.Range("E46") "  +0.04%  "

Set-TextValue $ws.Range("D47") "102.48"
Set-TextValue $ws.Range("E47") "  -0.04%  "

Set-TextValue $ws.Range("E48") "  -0.75%  "

Set-TextValue $ws.Range("D49") "1.766"
Set-TextValue $ws.Range("E49") "  -0.57%  "

Set-TextValue $ws.Range("D50") "9.160"
Set-TextValue $ws.Range("E50") "  -0.39%  "

Set-TextValue $ws.Range("D51") "7.017"
Set-TextValue $ws.Range("E51") "  +0.88%  "
